$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.157, -0.33,  -0.033, -116,  81,   42),
    @(0.103,  0.168,  -0.09,   157, -35, -168),
    @(0.117, -0.025,   0.472, -179,  68,  -71),
    @(0.303, -0.213,   0.294,   21,  46,  -69),
    @(-0.23, -0.194,   0.006, -135,  29,    7),
    @(0.247,  0.214,   0.504,   95,   5,  175),
    @(0.066, -0.212,   0.275,  -71, -63,  124),
    @(-0.107,-0.445,  -0.142,  -73, -42, -176),
    @(0.134,  0.232,   0.182,   74,  38,   56),
    @(-0.287,-0.418,   0.003, -166,   6, -109)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}

$wb.Save()
